$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value while forcing text storage (the source values are
# display strings like prices/percentages and must not be re-interpreted as numbers).
function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "26.718.64"
Set-TextValue $ws.Range("E2") "  +0.34%  "
Set-TextValue $ws.Range("D3") "1.602.02"
Set-TextValue $ws.Range("E3") "  +0.29%  "
Set-TextValue $ws.Range("E4") "  +0.31%  "
Set-TextValue $ws.Range("D5") "211.56"
Set-TextValue $ws.Range("E5") "  +0.09%  "
Set-TextValue $ws.Range("D6") "0.513"
Set-TextValue $ws.Range("E6") "  -0.43%  "
Set-TextValue $ws.Range("E7") "  +0.21%  "
Set-TextValue $ws.Range("E8") "  +0.32%  "
Set-TextValue $ws.Range("E9") "  +0.68%  "
Set-TextValue $ws.Range("D10") "19.75"
Set-TextValue $ws.Range("E10") "  +1.66%  "
Set-TextValue $ws.Range("E11") "  +0.57%  "
Set-TextValue $ws.Range("D12") "1.826.90"
Set-TextValue $ws.Range("E12") "  +0.28%  "
Set-TextValue $ws.Range("D13") "1.602.16"
Set-TextValue $ws.Range("E13") "  +0.80%  "
Set-TextValue $ws.Range("E14") "  +0.54%  "
Set-TextValue $ws.Range("E15") "  +0.60%  "
Set-TextValue $ws.Range("D16") "65.27"
Set-TextValue $ws.Range("E16") "  +0.78%  "
Set-TextValue $ws.Range("D17") "26.694.33"
Set-TextValue $ws.Range("E17") "  +0.25%  "
Set-TextValue $ws.Range("E18") "  +1.55%  "
Set-TextValue $ws.Range("D19") "7.23"
Set-TextValue $ws.Range("E19") "  +2.69%  "
Set-TextValue $ws.Range("D20") "210.54"
Set-TextValue $ws.Range("E20") "  +0.97%  "
Set-TextValue $ws.Range("E21") "  +0.24%  "
Set-TextValue $ws.Range("D22") "4.32"
Set-TextValue $ws.Range("E22") "  +1.32%  "
Set-TextValue $ws.Range("E23") "  +0.49%  "
Set-TextValue $ws.Range("D24") "8.99"
Set-TextValue $ws.Range("E24") "  +1.43%  "
Set-TextValue $ws.Range("D25") "143.48"
Set-TextValue $ws.Range("E25") "  -1.38%  "
Set-TextValue $ws.Range("D26") "1.01"
Set-TextValue $ws.Range("E26") "  +0.20%  "
Set-TextValue $ws.Range("E27") "  -0.32%  "
Set-TextValue $ws.Range("E28") "  -0.85%  "
Set-TextValue $ws.Range("D29") "15.42"
Set-TextValue $ws.Range("E29") "  +1.15%  "
Set-TextValue $ws.Range("E30") "  +1.72%  "
Set-TextValue $ws.Range("E31") "  -0.30%  "
Set-TextValue $ws.Range("D32") "3.27"
Set-TextValue $ws.Range("E32") "  +1.59%  "
Set-TextValue $ws.Range("E33") "  +1.81%  "
Set-TextValue $ws.Range("D34") "1.300.16"
Set-TextValue $ws.Range("E34") "  +2.44%  "
Set-TextValue $ws.Range("D36") "0.609"
Set-TextValue $ws.Range("E36") "  -2.22%  "
Set-TextValue $ws.Range("E37") "  +1.15%  "
Set-TextValue $ws.Range("E38") "  +22.48%  "
Set-TextValue $ws.Range("E39") "  -0.06%  "
Set-TextValue $ws.Range("E40") "  -1.83%  "
Set-TextValue $ws.Range("E41") "  -1.41%  "
Set-TextValue $ws.Range("E42") "  -0.14%  "
Set-TextValue $ws.Range("D43") "0.784"
Set-TextValue $ws.Range("E43") "  -0.28%  "
Set-TextValue $ws.Range("D44") "63.25"
Set-TextValue $ws.Range("E44") "  -1.42%  "
Set-TextValue $ws.Range("D45") "1.737.85"
Set-TextValue $ws.Range("E45") "  +0.12%  "
Set-TextValue $ws.Range("D46") "91.14"
Set-TextValue $ws.Range("E46") "  +1.35%  "
Set-TextValue $ws.Range("E47") "  -2.08%  "
Set-TextValue $ws.Range("D48") "0.0₆0105"
Set-TextValue $ws.Range("E48") "  -1.32%  "
Set-TextValue $ws.Range("E49") "  -0.88%  "
Set-TextValue $ws.Range("D50") "0.0518"
Set-TextValue $ws.Range("E50") "  +1.94%  "
Set-TextValue $ws.Range("E51") "  +0.01%  "
